# "Generate Report for Handback"
#
# This script fills in the previously-empty "Latest Target File" (E) and
# "Latest Handback File" (F) columns on the per-locale handoff/handback
# status sheets (zh-cn, de-de), updates the "Latest Handback DateTime" (G)
# for the real (non .localization-config) rows, and flips the Status text
# from "Ready for handoff" to "Handed back: in sync with en-US" everywhere
# it appears (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: just the status text changes (shared string, used by
# every sheet that references "Ready for handoff").
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusNew
$wsOverview.Range("C2").Value = $statusNew
$wsOverview.Range("B3").Value = $statusNew
$wsOverview.Range("C3").Value = $statusNew

# ---------------------------------------------------------------------
# Helper data describing the two locale sheets.
# ---------------------------------------------------------------------
$locales = @(
    @{
        SheetName = "zh-cn"
        MdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/0274075d79c8dc2146a263d9ffe1b03f03c9e711/e2e/bb390907-ecb5-45e7-b2e2-f0a195974103.md"
        XlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e5d05c6663e4ee562367a558a3717f41a69396f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bb390907-ecb5-45e7-b2e2-f0a195974103.ea38db53c51cec0207a32dfce0671018fd5962ab.zh-cn.xlf"
        MdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/0274075d79c8dc2146a263d9ffe1b03f03c9e711/e2e/fb906957-587a-453c-b3f2-1c11d9737d9e.md"
        XlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e5d05c6663e4ee562367a558a3717f41a69396f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/fb906957-587a-453c-b3f2-1c11d9737d9e.edb2dace96e23f75d446cd65ab7393c6f2b42b6a.zh-cn.xlf"
        HandbackDateTime = "2016-03-08 10:30:58"
    },
    @{
        SheetName = "de-de"
        MdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/0274075d79c8dc2146a263d9ffe1b03f03c9e711/e2e/bb390907-ecb5-45e7-b2e2-f0a195974103.md"
        XlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c3d3ccfddab948fd63bc84f7843c01ac7def7585/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bb390907-ecb5-45e7-b2e2-f0a195974103.ea38db53c51cec0207a32dfce0671018fd5962ab.de-de.xlf"
        MdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/0274075d79c8dc2146a263d9ffe1b03f03c9e711/e2e/fb906957-587a-453c-b3f2-1c11d9737d9e.md"
        XlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c3d3ccfddab948fd63bc84f7843c01ac7def7585/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/fb906957-587a-453c-b3f2-1c11d9737d9e.edb2dace96e23f75d446cd65ab7393c6f2b42b6a.de-de.xlf"
        HandbackDateTime = "2016-03-08 10:31:12"
    }
)

$mdDisplay1 = "bb390907-ecb5-45e7-b2e2-f0a195974103.md"
$mdDisplay2 = "fb906957-587a-453c-b3f2-1c11d9737d9e.md"
$xlfDisplay1Suffix = "ea38db53c51cec0207a32dfce0671018fd5962ab"
$xlfDisplay2Suffix = "edb2dace96e23f75d446cd65ab7393c6f2b42b6a"

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.SheetName)

    $xlfDisplay1 = "bb390907-ecb5-45e7-b2e2-f0a195974103." + $xlfDisplay1Suffix + "." + $locale.SheetName + ".xlf"
    $xlfDisplay2 = "fb906957-587a-453c-b3f2-1c11d9737d9e." + $xlfDisplay2Suffix + "." + $locale.SheetName + ".xlf"

    # Status column (B) -> handed back
    $ws.Range("B2").Value = $statusNew
    $ws.Range("B3").Value = $statusNew

    # Row 2 (bb390907...) : Latest Target File (E) + Latest Handback File (F)
    $ws.Hyperlinks.Add($ws.Range("E2"), $locale.MdUrl1, "", "", $mdDisplay1)
    $ws.Hyperlinks.Add($ws.Range("F2"), $locale.XlfUrl1, "", "", $xlfDisplay1)

    # Row 3 (fb906957...) : Latest Target File (E) + Latest Handback File (F)
    $ws.Hyperlinks.Add($ws.Range("E3"), $locale.MdUrl2, "", "", $mdDisplay2)
    $ws.Hyperlinks.Add($ws.Range("F3"), $locale.XlfUrl2, "", "", $xlfDisplay2)

    # Latest Handback DateTime (G) for the two real rows.
    $ws.Range("G2").Value = $locale.HandbackDateTime
    $ws.Range("G3").Value = $locale.HandbackDateTime
}
